# data/parameters.xlsx update
# 1. set concentration of species except VA, Pl, R1, and R2 (to 0)
# 2. change VA and Pl concentration to 2.22 µM (2.22E-6 M)
# (also tweaks several downstream kinetic-rate / receptor-density parameters
#  to their newly re-measured values, per the commit)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Species concentrations (rows 2-7, column C) ---
# VA (row 2) -> 2.22 uM
$ws.Range("C2").Value = 0.00000222
# VB (row 3) -> 0 (species other than VA, Pl, R1, R2)
$ws.Range("C3").Value = 0
# Pl (row 4) -> 2.22 uM
$ws.Range("C4").Value = 0.00000222
# PDAA (row 5) -> 0
$ws.Range("C5").Value = 0
# PDAB (row 6) -> 0
$ws.Range("C6").Value = 0
# PDBB (row 7) -> 0
$ws.Range("C7").Value = 0

# --- Receptor / co-receptor densities (rows 8-12, column C) ---
# R1 (row 8) stays a formula, value bumped 1600 -> 1604
$ws.Range("C8").Formula = "=1604"
# R2 (row 9) formula replaced by a plain literal value 4095
$ws.Range("C9").Value = 4095
# N1 (row 10) -> 0 (species other than VA, Pl, R1, R2)
$ws.Range("C10").Value = 0
# PDRa (row 11) -> 0
$ws.Range("C11").Value = 0
# PDRb (row 12) -> 0
$ws.Range("C12").Value = 0

# --- Updated kinetic-rate parameters ---
$ws.Range("C32").Value = 11370000
$ws.Range("C34").Value = 4650000
$ws.Range("C36").Value = 1260000
$ws.Range("C37").Value = 0.00348
$ws.Range("C43").Value = 0.00009145
$ws.Range("C55").Value = 0.0003035

$ws.Range("C64").Formula = "=(336000+325000)/2"
$ws.Range("C65").Formula = "=(0.000604+0.000905)/2"

# --- View state: selection grows to C5:C7 and the sheet is scrolled so
#     row 46 is the first visible row. Selection is supported by this
#     COM host; the scroll position is best-effort (no-op if unsupported).
$ws.Range("C5:C7").Select()
$excel.ActiveWindow.ScrollRow = 46
